# repull data, push all data, mean calculation
# Update the dSF (column F) values for a set of rows to reflect the
# repulled/recalculated data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 4
$ws.Range("F6").Value = -2
$ws.Range("F11").Value = -2
$ws.Range("F14").Value = 2
$ws.Range("F17").Value = -8
$ws.Range("F20").Value = -4
$ws.Range("F25").Value = -6
$ws.Range("F26").Value = -5
